$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$bleu  = "tube-alpex/Tube-Alpex-isolé-bleu.png"
$rouge = "tube-alpex/Tube-Alpex-isolé-rouge.png"
$plain = "tube-alpex/Tube-Alpex.png"

$ws.Range("A2").Value  = $bleu
$ws.Range("A3").Value  = $rouge
$ws.Range("A4").Value  = $bleu
$ws.Range("A5").Value  = $rouge
$ws.Range("A6").Value  = $bleu
$ws.Range("A7").Value  = $rouge
$ws.Range("A8").Value  = $plain
$ws.Range("A9").Value  = $plain
$ws.Range("A10").Value = $plain
$ws.Range("A11").Value = $plain
$ws.Range("A12").Value = $plain
$ws.Range("A13").Value = $bleu
$ws.Range("A14").Value = $rouge
$ws.Range("A15").Value = $bleu
$ws.Range("A16").Value = $plain
$ws.Range("A17").Value = $rouge
$ws.Range("A18").Value = $rouge
$ws.Range("A19").Value = $bleu

$ws.Range("A16").Select()
